$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.619706153869629
$ws.Range("B1").Value = 4.114301204681396
$ws.Range("C1").Value = 2.815733909606934
$ws.Range("D1").Value = 1.024818778038025
$ws.Range("E1").Value = 0.8866534829139709
